$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date-formatted style (s="4") from A75 down onto the new date cells
# so the new rows match the existing "Date" column formatting exactly.
$ws.Range("A75").Copy()
$ws.Range("A76:A79").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 76: 2010-04-04 (40272), 8 hours, "beagleboard avr-can interface"
$ws.Range("A76").Value = 40272
$ws.Range("B76").Value = 8
$ws.Range("C76").Value = "beagleboard avr-can interface"

# Row 77: 2010-04-05 (40273), 8 hours, "beagleboard avr-can interface"
$ws.Range("A77").Value = 40273
$ws.Range("B77").Value = 8
$ws.Range("C77").Value = "beagleboard avr-can interface"

# Row 78: 2010-04-06 (40274), 1.5 hours, "Group Meeting"
$ws.Range("A78").Value = 40274
$ws.Range("B78").Value = 1.5
$ws.Range("C78").Value = "Group Meeting"

# Row 79: 2010-04-06 (40274), 1 hour, "Weekly Meeting"
$ws.Range("A79").Value = 40274
$ws.Range("B79").Value = 1
$ws.Range("C79").Value = "Weekly Meeting"

# Update the view: scroll position and active selection, matching the
# author's final cursor position after entering the new rows.
$excel.ActiveWindow.ScrollRow = 51
$ws.Range("A80").Select()
